# Estadisticos Segundo Parcial 23 Mayo
# Update the "Estadisticos 2P" sheet with the 2nd-partial pass/fail stats,
# refresh the "Estadisticos Final" averages that depend on it, and drop the
# student who is no longer in need of a make-up exam from "Rescatables".

$wb = $excel.ActiveWorkbook

# --- Estadisticos 2P --------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

# Row 2: TEMAS DE FILOSOFIA / 6ALCM
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 1
$ws2P.Range("F2").Value = 40
$ws2P.Range("G2").Value = 97.56
$ws2P.Range("H2").Value = 8.199999999999999

# Row 3: TEMAS DE FILOSOFIA / 6BLCM
$ws2P.Range("D3").Value = 0
$ws2P.Range("E3").Value = 1
$ws2P.Range("F3").Value = 35
$ws2P.Range("G3").Value = 97.22
$ws2P.Range("H3").Value = 8.9

# Row 4: TEMAS DE FILOSOFIA / 6APV
$ws2P.Range("D4").Value = 0
$ws2P.Range("E4").Value = 0
$ws2P.Range("F4").Value = 17
$ws2P.Range("G4").Value = 100
$ws2P.Range("H4").Value = 8.800000000000001

# Row 5: TEMAS DE FILOSOFIA / 6ARHV
$ws2P.Range("D5").Value = 0
$ws2P.Range("E5").Value = 0
$ws2P.Range("F5").Value = 24
$ws2P.Range("G5").Value = 100
$ws2P.Range("H5").Value = 8.300000000000001

# --- Estadisticos Final ------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("H2").Value = 8.699999999999999
$wsFinal.Range("H3").Value = 9.1
$wsFinal.Range("H4").Value = 9

$wsFinal.Range("D5").Value = 0
$wsFinal.Range("E5").Value = 0
$wsFinal.Range("F5").Value = 24
$wsFinal.Range("G5").Value = 100
$wsFinal.Range("H5").Value = 8.800000000000001

# --- Rescatables ---------------------------------------------------------
# The student 22330051920205 (ANTONIO AGUILAR, EMILIANO GERARDO - 6ARHV)
# no longer needs to be rescued; remove that record, shifting the rest up.
$wsResc = $wb.Worksheets.Item("Rescatables")
$wsResc.Rows("2:2").Delete()
